# Regenerate save_data to use K (strikeouts) instead of Strike# in column G.
# This overwrites the existing literal values in column G (rows 2-36) with
# the newly computed K values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 2
    4  = 2
    5  = 3
    6  = 3
    7  = 1
    8  = 1
    9  = 5
    10 = 2
    11 = 1
    12 = 5
    13 = 4
    14 = 3
    15 = 7
    16 = 3
    17 = 2
    18 = 4
    19 = 4
    20 = 5
    21 = 8
    22 = 3
    23 = 2
    24 = 4
    25 = 6
    26 = 1
    27 = 1
    28 = 4
    29 = 2
    30 = 1
    31 = 2
    32 = 3
    33 = 4
    34 = 4
    35 = 3
    36 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
